$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "electrodes" column first (reset to [(0,0)] since a trial is
# now repeated if a stimulator error occurs), then the "channels" column.
$ws.Range("E2").Value = "[(0,0)]"
$ws.Range("E3").Value = "[(0,0)]"
$ws.Range("E4").Value = "[(0,0)]"
$ws.Range("E5").Value = "[(0,0)]"

$ws.Range("D2").Value = "[1, 2, 3]"
$ws.Range("D3").Value = "[1]"
$ws.Range("D4").Value = "[3]"
$ws.Range("D5").Value = "[1, 2]"

# Update the active selection to match the saved state in the workbook.
$ws.Range("D7").Select()
